$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.374.74'
$ws.Range("E2").Value = '  -0.15%  '
$ws.Range("D3").Value = '1.843.81'
$ws.Range("E3").Value = '  -0.29%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.61'
$ws.Range("E5").Value = '  +0.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6277'
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07483'
$ws.Range("E8").Value = '  -2.71%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2894'
$ws.Range("E9").Value = '  -0.92%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.36'
$ws.Range("E10").Value = '  -2.51%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07722'
$ws.Range("E11").Value = '  -0.32%  '
$ws.Range("D12").Value = '1.844.05'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.984'
$ws.Range("E13").Value = '  -1.10%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6786'
$ws.Range("E14").Value = '  -0.49%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001026'
$ws.Range("E15").Value = '  -5.06%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '82.13'
$ws.Range("E16").Value = '  -1.77%  '
$ws.Range("D17").Value = '2.106.82'
$ws.Range("E17").Value = '  -1.00%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.092'
$ws.Range("E18").Value = '  -1.54%  '
$ws.Range("D19").Value = '29.410.73'
$ws.Range("E19").Value = '  -0.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '228.61'
$ws.Range("E20").Value = '  -0.14%  '
$ws.Range("E21").Value = '  -1.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.366'
$ws.Range("E23").Value = '  -1.34%  '
$ws.Range("E24").Value = '  +0.23%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.94'
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.390'
$ws.Range("E27").Value = '  -0.33%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.53'
$ws.Range("E28").Value = '  -1.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.392'
$ws.Range("E29").Value = '  +2.77%  '
$ws.Range("E30").Value = '  +1.11%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.05686'
$ws.Range("E31").Value = '  +1.12%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.100'
$ws.Range("E32").Value = '  -0.62%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.039'
$ws.Range("E33").Value = '  -0.23%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.819'
$ws.Range("E34").Value = '  -1.37%  '
$ws.Range("E35").Value = '  -1.63%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6908'
$ws.Range("E36").Value = '  -2.39%  '
$ws.Range("E37").Value = '  -0.18%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.847'
$ws.Range("E38").Value = '  +3.36%  '
$ws.Range("D39").Value = '1.251.89'
$ws.Range("E39").Value = '  +2.22%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01814'
$ws.Range("E40").Value = '  +1.30%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.520'
$ws.Range("E41").Value = '  +1.04%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9042'
$ws.Range("E42").Value = '  +0.27%  '
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("D44").Value = '2.006.86'
$ws.Range("E44").Value = '  -1.31%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.22'
$ws.Range("E45").Value = '  -0.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '65.64'
$ws.Range("E46").Value = '  -0.65%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.078'
$ws.Range("E47").Value = '  -1.55%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1161'
$ws.Range("E48").Value = '  +0.46%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.965'
$ws.Range("E49").Value = '  -0.35%  '
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00000000114'
$ws.Range("E50").Value = '  -4.64%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3938'
$ws.Range("E51").Value = '  -2.08%  '
